# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-row suffixes to "_FV2310" / "_FV2404"
# - Turn the data range into an Excel Table ("Table1")
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header row (row 1) cells: *_old -> *_FV2310, *_new -> *_FV2404
# ---------------------------------------------------------------------------
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J -> "<name>_FV2310"
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($oldHeaders[$i] + "_FV2310")
}

# Column K stays "diff" (unchanged)
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L..U -> "<name>_FV2404"
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($oldHeaders[$i] + "_FV2404")
}

# ---------------------------------------------------------------------------
# 2. Convert the used range into a table
# ---------------------------------------------------------------------------
$range = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $range, $true)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3. Freeze the header row
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
